$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.510.36'
$ws.Range("E2").Value = '  +0.13%  '

# Row 3
$ws.Range("D3").Value = '2.105.67'
$ws.Range("E3").Value = '  +0.53%  '

# Row 4
$ws.Range("E4").Value = '  +0.70%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '335.64'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.71%  '

# Row 6
$ws.Range("E6").Value = '  +0.64%  '

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.5219'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.19%  '

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.4528'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +3.61%  '

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '55.28'
$c.Style = "Normal"

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.09023'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +2.07%  '

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '1.168'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +1.26%  '

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '24.57'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +1.06%  '

# Row 13
$ws.Range("D13").Value = '2.112.90'
$ws.Range("E13").Value = '  +1.08%  '

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '6.817'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +2.18%  '

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '8.087'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +5.58%  '

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.00001166'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +4.09%  '

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '96.90'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +1.14%  '

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '1.010'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.67%  '

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.06680'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +1.39%  '

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '19.36'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.56%  '

# Row 21
$ws.Range("E21").Value = '  +0.67%  '

# Row 22
$ws.Range("E22").Value = '  -0.43%  '

# Row 23
$ws.Range("D23").Value = '30.576.33'
$ws.Range("E23").Value = '  +0.19%  '

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '12.82'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +4.73%  '

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.356'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.71%  '

# Row 26
$ws.Range("D26").Value = '2.351.78'
$ws.Range("E26").Value = '  +0.68%  '

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '22.24'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.03%  '

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '163.22'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.54%  '

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '2.505'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -2.01%  '

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '133.27'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.42%  '

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.210'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +2.24%  '

# Row 32
$ws.Range("E32").Value = '  +0.05%  '

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '1.638'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.48%  '

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '6.334'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +3.25%  '

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '3.953'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +1.25%  '

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '10.39'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +2.26%  '

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '5.891'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +8.48%  '

# Row 38
$ws.Range("E38").Value = '  +1.63%  '

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.06797'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.34%  '

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.2312'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +2.66%  '

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '12.58'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -0.74%  '

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.6839'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -0.55%  '

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.255'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -0.59%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '14.21'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +1.76%  '

# Row 45
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.6429'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +1.19%  '

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '2.302'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +4.95%  '

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '3.671'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +1.38%  '

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.00000000354'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +18.57%  '

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.247'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.60%  '

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.3399'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +14.35%  '

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '83.04'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +1.71%  '
